$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: the Price column (D) stores numeric-looking text such as "246.00",
# "9.10", "0.0780", "2.80" where the trailing zero(s) matter. Excel's COM
# layer auto-converts a numeric-looking string assigned to .Value into a
# real number (dropping trailing zeros), so each such D cell is switched to
# Text format ("@") immediately before its value is written. (Multi-area
# comma ranges don't reliably propagate NumberFormat here, so each cell is
# set individually.)

# Row 2 - Bitcoin
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "47.212.12"
$ws.Range("E2").Value = "  -0.85%  "

# Row 3 - Ethereum
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.484.64"
$ws.Range("E3").Value = "  -0.57%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.00%  "

# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "320.86"
$ws.Range("E5").Value = "  -0.92%  "

# Row 6 - Solana
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "108.01"
$ws.Range("E6").Value = "  +2.64%  "

# Row 7 - XRP
$ws.Range("E7").Value = "  -0.32%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.01%  "

# Row 9 - Cardano
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.538"
$ws.Range("E9").Value = "  -0.21%  "

# Row 10 - Avalanche
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.07"
$ws.Range("E10").Value = "  +3.49%  "

# Row 11 - Dogecoin
$ws.Range("E11").Value = "  -0.49%  "

# Row 12 - TRON
$ws.Range("E12").Value = "  +0.60%  "

# Row 13 - Chainlink
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.42"
$ws.Range("E13").Value = "  +0.75%  "

# Row 14 - Polkadot
$ws.Range("E14").Value = "  +0.34%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.872.34"
$ws.Range("E15").Value = "  -0.53%  "

# Row 16 - WrappedEther
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.484.44"
$ws.Range("E16").Value = "  -0.72%  "

# Row 17 - Polygon
$ws.Range("E17").Value = "  +0.01%  "

# Row 18 - WrappedBTC
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "47.129.11"
$ws.Range("E18").Value = "  -0.73%  "

# Row 19 - InternetComputer(DFINITY)
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.49"
$ws.Range("E19").Value = "  +6.52%  "

# Row 20 - Uniswap
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.61"
$ws.Range("E20").Value = "  +0.93%  "

# Row 21 <-> Row 22 swap: ShibaInu/ImmutableX order flip, with updated values
$ws.Range("B21").Value = "ImmutableX"
$ws.Range("C21").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.80"
$ws.Range("E21").Value = "  +17.01%  "

$ws.Range("B22").Value = "ShibaInu"
$ws.Range("C22").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0₃0942"
$ws.Range("E22").Value = "  +0.87%  "

# Row 23 - Litecoin
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.48"
$ws.Range("E23").Value = "  -0.45%  "

# Row 24 - BitcoinCash
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "246.00"
$ws.Range("E24").Value = "  -2.06%  "

# Row 25 - PancakeSwap
$ws.Range("E25").Value = "  -1.28%  "

# Row 26 - Dai
$ws.Range("E26").Value = "  -0.05%  "

# Row 27 - EthereumClassic
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.63"
$ws.Range("E27").Value = "  -2.32%  "

# Row 28 - Toncoin
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.28"
$ws.Range("E28").Value = "  +3.64%  "

# Row 29 - Cosmos
$ws.Range("E29").Value = "  -1.51%  "

# Row 30 - Kaspa
$ws.Range("E30").Value = "  +4.24%  "

# Row 31 - InjectiveProtocol
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "34.57"
$ws.Range("E31").Value = "  -2.02%  "

# Row 32 - OKB
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "49.78"
$ws.Range("E32").Value = "  +0.57%  "

# Row 33 - Celestia
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.32"
$ws.Range("E33").Value = "  +1.67%  "

# Row 34 - Filecoin
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.32"
$ws.Range("E34").Value = "  -1.16%  "

# Row 35 - Hedera
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0780"
$ws.Range("E35").Value = "  -0.27%  "

# Row 36 - FirstDigitalUSD
$ws.Range("E36").Value = "  +0.16%  "

# Row 37 - RenderToken
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.79"
$ws.Range("E37").Value = "  +3.31%  "

# Row 38 - ARBITRUM
$ws.Range("E38").Value = "  -0.03%  "

# Row 39 - LidoDAOToken
$ws.Range("E39").Value = "  -2.02%  "

# Row 40 - EnergySwap
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "22.92"
$ws.Range("E40").Value = "  +8.67%  "

# Row 41 - Stellar
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.111"
$ws.Range("E41").Value = "  -0.04%  "

# Row 42 <-> Row 43 swap: WEMIXToken/Monero order flip, with updated values
$ws.Range("B42").Value = "Monero"
$ws.Range("C42").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "120.93"
$ws.Range("E42").Value = "  -0.17%  "

$ws.Range("B43").Value = "WEMIXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.21"
$ws.Range("E43").Value = "  -1.04%  "

# Row 44 - VeChain
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0297"
$ws.Range("E44").Value = "  +0.10%  "

# Row 45 - Maker
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.995.63"
$ws.Range("E45").Value = "  +1.36%  "

# Row 46 - NEARProtocol
$ws.Range("E46").Value = "  +1.23%  "

# Row 47 - ApeXProtocol
$ws.Range("E47").Value = "  -4.46%  "

# Row 48 - Stacks
$ws.Range("E48").Value = "  -2.89%  "

# Row 49 - FraxShare
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.10"
$ws.Range("E49").Value = "  -1.06%  "

# Row 50 - THORChain
$ws.Range("E50").Value = "  -4.60%  "

# Row 51 - MultiversX
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "56.81"
$ws.Range("E51").Value = "  +3.42%  "
